$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = "RJ TV 1"
$ws.Range("C2").Value = "Defesa Civil"
$ws.Range("D2").Value = "2025-04-02T12:24"
$ws.Range("F2").Value = "Alerta no RJ. Estado deve ter chuva muito forte a partir de sexta-feira. Repórter *ao vivo* do Rio de Janeiro. "

# Row 3 updates
$ws.Range("B3").Value = "RJ TV 1"
$ws.Range("C3").Value = "PROCON"
$ws.Range("D3").Value = "2025-04-02T12:51"
$ws.Range("E3").Value = "Positivo"
$ws.Range("F3").Value = "Variação de preços dos produtos da Páscoa. Em Campos, a fiscalização encontrou 75% de variação no preço de alguns produtos. Repórter *ao vivo*."

# Remove rows 4, 5, 6
$ws.Rows("4:6").Delete()
